$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.252.46'
$ws.Range("E2").Value = '  +1.99%  '
$ws.Range("D3").Value = '1.814.64'
$ws.Range("E3").Value = '  +3.49%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4356'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.64%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3667'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.92'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.96%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07674'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.143'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.08%  '
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.04'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.318'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.488'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.23%  '
$ws.Range("D16").Value = '1.825.64'
$ws.Range("E16").Value = '  +4.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '95.06'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +7.59%  '
$ws.Range("E18").Value = '  +1.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06478'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9999'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.41'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.81%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.233'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.28%  '
$ws.Range("D23").Value = '28.272.85'
$ws.Range("E23").Value = '  +2.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.58'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.130'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.30'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.93%  '
$ws.Range("E27").Value = '  +0.79%  '
$ws.Range("D28").Value = '2.025.68'
$ws.Range("E28").Value = '  +3.92%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.278'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '130.84'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.211'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.57%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.032'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09134'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("E34").Value = '  -2.37%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.08'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02387'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.241'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2177'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6593'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06206'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.201'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.49%  '
$ws.Range("E42").Value = '  +1.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.428'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9993'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.84'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6104'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.739'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.28%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.79'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.08%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.023'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.46%  '
$ws.Range("E50").Value = '  +3.93%  '
$ws.Range("E51").Value = '  +1.46%  '
